# Populate the header row for the first-year students sheet with
# USN / Name / Age / Phone columns (A1:D1), as introduced by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"

# Leave the cursor on the last entered cell, matching the saved selection.
$ws.Range("D1").Select()
